$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.979.98'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.760.48'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.57%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '335.01'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.20%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3936'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.56%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3392'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.60%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.34'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.118'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07218'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.31'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -4.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.140'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -5.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.096'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.695.03'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -5.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001058'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.89%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06625'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '80.38'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.41%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.94'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.225'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -3.72%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.974.67'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.14%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.387'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.52'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.93'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -4.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.312'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.97%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.920.29'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -3.51%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.276'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -11.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '129.34'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.62%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.810'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08725'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.06'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06180'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02289'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -5.91%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.137'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6487'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -5.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2112'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.497'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.21%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.204'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -3.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9993'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.878'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -5.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.79'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.829'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5983'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -5.50%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '126.48'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -5.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.986'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -4.82%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.157'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07004'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -6.58%  '
